$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds account numbers that must stay text even though they are
# all-digits (otherwise Excel coerces them to numbers / scientific notation).
$ws.Range("C2:C6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "MOHAMED BADRANE"
$ws.Range("B2").Value = "I83603"
$ws.Range("C2").Value = "225400000805987601012173"
$ws.Range("D2").Value = "KHOURIBGA"
$ws.Range("E2").Value = "CA"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "605/KHOURIBGA NAHDA"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 7500
$ws.Range("J2").Value = 375
$ws.Range("K2").Value = 7125

# Row 3
$ws.Range("A3").Value = "ZERNAKH ABDELLAH"
$ws.Range("B3").Value = "IB19558"
$ws.Range("C3").Value = "145101211406073828000084"
$ws.Range("D3").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Point de vente"
$ws.Range("G3").Value = "052/FKIH BEN SALEH"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 11000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 11000

# Row 4
$ws.Range("A4").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B4").Value = "Q251990"
$ws.Range("C4").Value = "007400000313200019604463"
$ws.Range("D4").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E4").Value = "AWB"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 6750
$ws.Range("J4").Value = 675
$ws.Range("K4").Value = 6075

# Row 5
$ws.Range("A5").Value = "NOUBAIL MOHAMMED"
$ws.Range("B5").Value = "IR801997"
$ws.Range("C5").Value = "007400000313200019604463"
$ws.Range("D5").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E5").Value = "AWB"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 6750
$ws.Range("J5").Value = 675
$ws.Range("K5").Value = 6075

# Row 6 (totals row, blank labels, totals numbers)
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 32000
$ws.Range("J6").Value = 1725
$ws.Range("K6").Value = 30275

# Remove old rows 7 and 8 (delete entire rows, shifting cells up)
$ws.Range("A7:K8").Delete()
